# edit.ps1 - applies the CatSwarmRapordocx commit:
#  1) Drop the stray combining-cedilla (U+00B8) typo after "s" in two words
#  2) Expand the "Adim 6:" terminal paragraph into full sentence + a new
#     "Kedi Suerues Optimizasyon Algoritmasi" sub-section (16 new paragraphs)
#  3) Move the hidden _GoBack bookmark from the title block down to the
#     blank paragraph that now follows the new section (tracks last edit)

$d = $word.ActiveDocument

# --- 1a. davranıs¸ının -> davranısının (single occurrence in the doc) ---
$null = $d.Content.Find.Execute("davranıs¸ının", $false, $false, $false, $false, $false, $true, 1, $false, "davranısının", 2)

# --- 1b. çıkarılmıs¸tır. -> çıkarılmıstır. (only the SECOND occurrence changes;
#         the first one, followed by "Kedilerin", must stay untouched) ---
$markerRange = $d.Content
$null = $markerRange.Find.Execute("çıkarılmıs¸tır. Kedilerin", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterFirstHit = $markerRange.End
$secondRange = $d.Range($afterFirstHit, $d.Content.End)
$null = $secondRange.Find.Execute("çıkarılmıs¸tır. ", $false, $false, $false, $false, $false, $true, 1, $false, "çıkarılmıstır. ", 2)

# --- 2. Replace the "Adım 6:" paragraph with the expanded text plus the
#        new "Kedi Sürüsü Optimizasyon Algoritması" sub-section that follows it ---
$findStep6 = $d.Content
$null = $findStep6.Find.Execute("Adım 6:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$step6Para = $findStep6.Paragraphs(1)

$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>Adım 6</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Sonlandırma (bitirme) koşulları sağlanmışsa programı sonlandır, aksi durumda Adım 3’ten Adım 5’e kadar tekrar et</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Kedi Sürüsü Optimizasyon Algoritması</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">KSO’nun, arama modu ve izleme modu adında iki alt modu </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">vardır. Bu iki modu algoritma şeklinde birleştirmek için, </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">arama moduyla izleme modunu birleştirmeyi sağlayan bir </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">karışım oranı (KO) tanımlanmaktadır. Kediler dinlenme </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">zamanında hareket etmeye karar verdiklerinde, hareket çok </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">dikkatli ve yavaşça yapılmaktadır. Bu hareket, arama </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">moduna yansıtılmaktadır. İzleme modu kedi tarafından bir </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">hedefin takip edilmesini modellemektedir. Kediler, enerji </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">kaynaklarını fazla kullanmalarına yol açan objeleri takip </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">etmeye çok az zaman harcamaktadır. Kedilerin </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">zamanlarının çoğunu dinlenmeye ve gözetlemeye (mesela </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>zamanlarının çoğu arama modunda geçmektedir)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">harcadığını garantilemek için KO’ya çok küçük bir değer </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">atanmaktadır. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:bidi w:val="0"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:lang w:val="tr-TR" w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>

'@

$step6Range = $step6Para.Range
$step6Range.Collapse(0)
$null = $step6Range.InsertXML($newBlockXml)

# The paragraph immediately after the (now much larger) inserted block is
# the blank paragraph that used to sit right after "Adım 6:" (sz=34 spacer).
$bookmarkDestPara = $step6Para.Next()

# --- 3. Move the hidden _GoBack bookmark to that blank paragraph ---
$d.Bookmarks.Add("_GoBack", $bookmarkDestPara.Range)

Write-Output "edit applied"
